$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.360.09"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.711.21"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5295"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2668"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06665"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07693"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.519"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "1.945.65"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "1.689.66"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5830"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "0.0₅8203"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "27.350.71"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "223.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.642"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.63%  "
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.265"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.297"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.461"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.437"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.869"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01640"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").Value = "1.081.11"
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.814"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8455"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "1.853.04"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4535"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "

Write-Output "Applied cryptos list update"
